$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the overall budget for home appliances (column H, row 2) from 500000 to 200000
$ws.Range("H2").Value = 200000

# Move the active cell from H3 to H2 (as reflected by the sheet's recorded selection)
$ws.Range("H2").Select()

# Touch the font on the "Mobile Number" header (C1) so Excel (re)writes its font
# metadata, matching the style update recorded for that cell
$ws.Range("C1").Font.Name = "Arial"
